$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price figures stored as literal text in the source sheet
# (dotted thousand separators like "3.105.78" are not valid numbers anyway).
# A leading apostrophe forces Excel to keep/treat the assignment as text,
# same as a manual text entry, without touching the cells number format.

$ws.Range("D2").Value = "'62.805.05"
$ws.Range("D3").Value = "'3.105.78"
$ws.Range("E3").Value = '  +2.61%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'585.32"
$ws.Range("E5").Value = '  +3.49%  '
$ws.Range("D6").Value = "'143.66"
$ws.Range("E6").Value = '  +2.24%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = "'3.096.74"
$ws.Range("E8").Value = '  +2.66%  '
$ws.Range("E9").Value = '  +1.43%  '
$ws.Range("D10").Value = "'0.149"
$ws.Range("E10").Value = '  +11.14%  '
$ws.Range("D11").Value = "'5.70"
$ws.Range("E12").Value = '  +1.10%  '
$ws.Range("D13").Value = "'0.0000245"
$ws.Range("E13").Value = '  +5.12%  '
$ws.Range("D14").Value = "'35.36"
$ws.Range("E14").Value = '  +3.65%  '
$ws.Range("E15").Value = '  -0.04%  '
$ws.Range("D16").Value = "'3.619.01"
$ws.Range("E16").Value = '  +2.62%  '
$ws.Range("D17").Value = "'7.21"
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").Value = "'3.103.42"
$ws.Range("E18").Value = '  +2.63%  '
$ws.Range("D19").Value = "'62.742.75"
$ws.Range("E19").Value = '  +4.94%  '
$ws.Range("E20").Value = '  +5.90%  '
$ws.Range("D21").Value = "'14.05"
$ws.Range("E21").Value = '  +2.58%  '
$ws.Range("D22").Value = "'0.728"
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("D23").Value = "'7.51"
$ws.Range("E23").Value = '  +5.06%  '
$ws.Range("D24").Value = "'13.40"
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("D25").Value = "'82.17"
$ws.Range("E25").Value = '  +1.44%  '
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("E27").Value = '  -1.15%  '
$ws.Range("D28").Value = "'2.68"
$ws.Range("E28").Value = '  +4.76%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = "'8.25"
$ws.Range("E30").Value = '  +4.88%  '
$ws.Range("D31").Value = "'6.80"
$ws.Range("E31").Value = '  +7.71%  '
$ws.Range("D32").Value = "'26.91"
$ws.Range("E32").Value = '  +3.23%  '
$ws.Range("E33").Value = '  +8.48%  '
$ws.Range("D34").Value = "'0.0₃0828"
$ws.Range("E34").Value = '  +4.94%  '
$ws.Range("D35").Value = "'2.35"
$ws.Range("E35").Value = '  +10.50%  '
$ws.Range("E36").Value = '  +3.07%  '
$ws.Range("D37").Value = "'6.02"
$ws.Range("E37").Value = '  +0.88%  '
$ws.Range("D38").Value = "'3.15"
$ws.Range("E38").Value = '  +11.62%  '
$ws.Range("D39").Value = "'50.96"
$ws.Range("E39").Value = '  +3.57%  '
$ws.Range("D40").Value = "'8.77"
$ws.Range("E40").Value = '  +0.92%  '
$ws.Range("D41").Value = "'429.52"
$ws.Range("E41").Value = '  +6.26%  '
$ws.Range("D42").Value = "'2.905.26"
$ws.Range("E42").Value = '  +4.05%  '
$ws.Range("E43").Value = '  +3.47%  '
$ws.Range("D44").Value = "'0.277"
$ws.Range("E44").Value = '  +8.26%  '
$ws.Range("E45").Value = '  +2.65%  '
$ws.Range("D46").Value = "'2.16"
$ws.Range("E46").Value = '  +6.53%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = "'123.92"
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("B49").Value = 'Arweave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D49").Value = "'34.89"
$ws.Range("E49").Value = '  +2.77%  '
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("D51").Value = "'24.69"
$ws.Range("E51").Value = '  +4.20%  '
